# Update the mobiles worksheet: rename headers, add a "Reviews" column (D),
# and refresh product rows 2-21 with the latest scraped listing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename existing headers and add new "Reviews" header ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Price"
$ws.Range("C1").Value = "Rating"

# New column D needs the same bold/centered/bordered header style as A1:C1,
# so copy formatting from the existing header cell C1 onto D1.
$ws.Range("D1").Value = "Reviews"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Reviews"
$excel.CutCopyMode = $false

# Price/Rating/Reviews columns hold numeric-looking text (e.g. "4.5", "5,050")
# that must stay text, not be auto-converted to numbers by Excel.
$ws.Range("B2:D21").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("A2").Value = "iPhone Air 256 GB: Thinnest iPhone Ever, 16.63 cm (6.5″) Display with Promotion up to 120Hz, Powerful A19 Pro Chip, Center Stage Front Camera, All-Day Battery Life; Light Gold"
$ws.Range("B2").Value = "1,19,900"
$ws.Range("C2").Value = "4.5"
$ws.Range("D2").Value = "5,050"

# --- Row 3 ---
$ws.Range("A3").Value = "Apple iPhone 15 (128 GB) - Blue"
$ws.Range("B3").Value = "47,999"
$ws.Range("C3").Value = "4.0"
$ws.Range("D3").Value = "6,910"

# --- Row 4 ---
$ws.Range("A4").Value = "Samsung Galaxy M05 (Mint Green, 4GB RAM, 64 GB Storage) | 50MP Dual Camera | Bigger 6.7`" HD+ Display | 5000mAh Battery | 25W Fast Charging | 2 Gen OS Upgrade & 4 Year Security Update | Without Charger"
$ws.Range("B4").Value = "6,249"
$ws.Range("C4").Value = "4.0"
$ws.Range("D4").Value = "2,160"

# --- Row 5 ---
$ws.Range("A5").Value = "iQOO Z10 Lite 5G (Cyber Green, 6GB RAM, 128GB Storage) | 6000 mAh Battery | Dimensity 6300 5G Processor with 433K+* AnTuTu Score | IP64 Rated & Military Grade Shock-Resistance*"
$ws.Range("B5").Value = "10,998"
$ws.Range("C5").Value = "4.8"
$ws.Range("D5").Value = "5"

# --- Row 6 ---
$ws.Range("A6").Value = "Pova Slim 5G (Slim White, 8+128GB) | World's Slimmest and Lightest* 5G with 5160mAh Battery | World's 1st Dynamic Mood Light | 144Hz 1.5K 3D Curve AMOLED | Military Grade MIL 810H Protection | IP64"
$ws.Range("B6").Value = "19,999"
$ws.Range("C6").Value = "3.9"
$ws.Range("D6").Value = "2,939"

# --- Row 7 ---
$ws.Range("A7").Value = "Samsung Galaxy M06 5G (Sage Green, 6GB RAM, 128 GB Storage) | MediaTek Dimensity 6300 | AnTuTu Score 422K+ | 12 5G Bands| 25W Fast Charging | 4 Gen of OS Upgrades | Without Charger"
$ws.Range("B7").Value = "8,999"
$ws.Range("C7").Value = "3.9"
$ws.Range("D7").Value = "2,271"

# --- Row 8 ---
$ws.Range("A8").Value = "Samsung Galaxy M16 5G (Thunder Black, 6GB RAM, 128 GB Storage) | MediaTek Dimensity 6300 | AnTuTu Score 422K+ | Super Amoled Display | AI | 25W Fast Charging | 6 Gen of OS Upgrades | Without Charger"
$ws.Range("B8").Value = "11,749"
$ws.Range("C8").Value = "4.0"
$ws.Range("D8").Value = "556"

# --- Row 9 ---
$ws.Range("A9").Value = "POCO C71, Cool Blue (6GB, 128GB)"
$ws.Range("B9").Value = "6,798"
$ws.Range("C9").Value = "4.2"
$ws.Range("D9").Value = "9,493"

# --- Row 10 ---
$ws.Range("A10").Value = "realme 14X 5G Smartphone Dimensity 6300 6nm Chip 8GB RAM 128GB ROM 6.67 Inch HD+ IP69 Waterproof 6000mAh Battery 45W Fast Charge Support Fingerprint GPS WiFi Bluetooth (Gold)"
$ws.Range("B10").Value = "14,190"
$ws.Range("C10").Value = "4.5"
$ws.Range("D10").Value = "1,464"

# --- Row 11 ---
$ws.Range("A11").Value = "OnePlus Nord CE4 Lite 5G (Super Silver, 8GB RAM, 128GB Storage) | Lifetime Display Warranty | 5500 mAh Battery, 80W SUPERVOOC and Reverse Charging | 50MP Camera with OIS | 120Hz AMOLED Display"
$ws.Range("B11").Value = "16,998"
$ws.Range("C11").Value = "4.6"
$ws.Range("D11").Value = "1,794"

# --- Row 12 ---
$ws.Range("A12").Value = "iPhone 16 128 GB: 5G Mobile Phone with Camera Control, A18 Chip and a Big Boost in Battery Life. Works with AirPods; White"
$ws.Range("B12").Value = "69,499"
$ws.Range("C12").Value = "4.0"
$ws.Range("D12").Value = "2,016"

# --- Row 13 ---
$ws.Range("A13").Value = "Apple iPhone 15 Plus (128 GB) - Black"
$ws.Range("B13").Value = "68,999"
$ws.Range("C13").Value = "3.9"
$ws.Range("D13").Value = "2,271"

# --- Row 14 ---
$ws.Range("A14").Value = "Redmi 13 5G Prime Edition, Orchid Pink, 8GB+128GB | India Debut SD 4 Gen 2 AE | 108MP Pro Grade Camera | 6.79in Largest Display in Segment"
$ws.Range("B14").Value = "11,199"
$ws.Range("C14").Value = "4.4"
$ws.Range("D14").Value = "1,138"

# --- Row 15 ---
$ws.Range("A15").Value = "Samsung Galaxy M16 5G (Blush Pink, 4GB RAM, 128 GB Storage) | MediaTek Dimensity 6300 | AnTuTu Score 422K+ | Super Amoled Display | AI | 25W Fast Charging | 6 Gen of OS Upgrades | Without Charger"
$ws.Range("B15").Value = "10,499"
$ws.Range("C15").Value = "4.0"
$ws.Range("D15").Value = "556"

# --- Row 16 ---
$ws.Range("A16").Value = "realme NARZO 80 Pro 5G (Speed Silver,12GB+256GB) | Segment's 1st MediaTek Dimensity 7400 Chipset | 6000mAh Titan Battery + 80W Ultra Charge | 4500nits HyperGlow Esports Display | IP69 Waterproof"
$ws.Range("B16").Value = "21,498"
$ws.Range("C16").Value = "3.9"
$ws.Range("D16").Value = "103"

# --- Row 17 ---
$ws.Range("A17").Value = "POCO C71, Desert Gold (6GB, 128GB)"
$ws.Range("B17").Value = "6,799"
$ws.Range("C17").Value = "3.8"
$ws.Range("D17").Value = "19,655"

# --- Row 18 ---
$ws.Range("A18").Value = "Redmi 15 5G Midnight Black 8GB + 256GB | Segment's Largest 7000mAhA Battery | Segment's Largest Display 17.53cm(6.9) Up to 144Hz | Snapdragon 6s Gen 3 | 18W Reverse Charging | 50MP AI Dual Camera"
$ws.Range("B18").Value = "16,998"
$ws.Range("C18").Value = "4.0"
$ws.Range("D18").Value = "2,016"

# --- Row 19 ---
$ws.Range("A19").Value = "Nokia 105 Classic | Single SIM Keypad Phone with Built-in UPI Payments, Long-Lasting Battery, Wireless FM Radio, Without Charger | Charcoal"
$ws.Range("B19").Value = "949"
$ws.Range("C19").Value = "3.9"
$ws.Range("D19").Value = "2,939"

# --- Row 20 ---
$ws.Range("A20").Value = "Redmi 13 5G Prime Edition, Black Diamond, 8GB+128GB | India Debut SD 4 Gen 2 AE | 108MP Pro Grade Camera | 6.79in Largest Display in Segment"
$ws.Range("B20").Value = "11,199"
$ws.Range("C20").Value = "4.3"
$ws.Range("D20").Value = "2,446"

# --- Row 21 ---
$ws.Range("A21").Value = "Samsung Galaxy M06 5G (Sage Green, 4GB RAM, 64GB Storage) | MediaTek Dimensity 6300 | AnTuTu Score 422K+ | 12 5G Bands| 25W Fast Charging | 4 Gen of OS Upgrades | Without Charger"
$ws.Range("B21").Value = "7,499"
$ws.Range("C21").Value = "3.9"
$ws.Range("D21").Value = "103"

Write-Host "Workbook updated: mobiles data refreshed with Reviews column."